$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet (Plan1 -> Sprint 2) ---
$ws.Name = "Sprint 2"

# --- Insert a new column D ("HORAS TRABALHADAS") before the existing STATUS column,
#     which shifts STATUS from D to E ---
$ws.Columns("D:D").Insert()

# --- Resize the table (Tabela2) so it covers the new column too ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E4"))

# --- Header + data for the new "HORAS TRABALHADAS" column ---
$ws.Range("D1").Value = "HORAS TRABALHADAS"
$ws.Range("D2").Value = "2 Horas 30 Minutos"
$ws.Range("D3").Value = "5 Horas"
$ws.Range("D4").Value = "1 Hora e 30 Minutos"

# --- Restore the STATUS header text in its new column position ---
$ws.Range("E1").Value = "STATUS"

# --- Give the new column's data the same centered look as the rest of the table,
#     plus a thin border around each cell ---
$dataRng = $ws.Range("D2:D4")
$dataRng.HorizontalAlignment = -4108
$dataRng.VerticalAlignment = -4108
$dataRng.Borders.LineStyle = 1
$dataRng.Borders.Weight = 2

# --- Table visual style ---
$lo.TableStyle = "TableStyleMedium6"

# --- Column widths (new HORAS TRABALHADAS column wider, STATUS keeps its old width) ---
$ws.Columns("D:D").ColumnWidth = 24.5
$ws.Columns("E:E").ColumnWidth = 11.5

# --- Sheet view: hide gridlines and move the active selection ---
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("F4").Select()
